$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be read/written as text so numeric-looking
# strings (e.g. "582.84", "1.00") are preserved verbatim, matching
# the original inlineStr cell contents instead of being coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '69.701.57'
$ws.Range("E2").Value = '  +0.69%  '
$ws.Range("D3").Value = '3.424.82'
$ws.Range("E3").Value = '  +1.12%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '582.84'
$ws.Range("E5").Value = '  -0.81%  '
$ws.Range("D6").Value = '176.55'
$ws.Range("E6").Value = '  -2.06%  '
$ws.Range("D7").Value = '3.416.54'
$ws.Range("E7").Value = '  +1.09%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("E10").Value = '  +2.94%  '
$ws.Range("D11").Value = '0.583'
$ws.Range("E11").Value = '  -1.24%  '
$ws.Range("D12").Value = '48.94'
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("D13").Value = '0.0000282'
$ws.Range("E13").Value = '  +0.40%  '
$ws.Range("D14").Value = '690.34'
$ws.Range("E14").Value = '  +2.01%  '
$ws.Range("D15").Value = '3.970.33'
$ws.Range("E15").Value = '  +1.01%  '
$ws.Range("E16").Value = '  -0.26%  '
$ws.Range("D17").Value = '69.730.07'
$ws.Range("E17").Value = '  +0.70%  '
$ws.Range("D18").Value = '3.424.86'
$ws.Range("E18").Value = '  +1.24%  '
$ws.Range("E19").Value = '  +0.94%  '
$ws.Range("D20").Value = '17.66'
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("D21").Value = '11.38'
$ws.Range("E21").Value = '  +0.40%  '
$ws.Range("D22").Value = '0.895'
$ws.Range("E22").Value = '  -0.78%  '
$ws.Range("D23").Value = '5.50'
$ws.Range("E23").Value = '  +1.44%  '
$ws.Range("D24").Value = '16.90'
$ws.Range("E24").Value = '  -1.19%  '
$ws.Range("D25").Value = '100.61'
$ws.Range("E25").Value = '  -2.66%  '
$ws.Range("D26").Value = '3.91'
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("E27").Value = '  -2.61%  '
$ws.Range("D28").Value = '9.56'
$ws.Range("E28").Value = '  -0.38%  '
$ws.Range("D29").Value = '33.43'
$ws.Range("E29").Value = '  -2.13%  '
$ws.Range("D30").Value = '8.72'
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("D31").Value = '7.12'
$ws.Range("E31").Value = '  +1.71%  '
$ws.Range("D32").Value = '572.75'
$ws.Range("E32").Value = '  +3.23%  '
$ws.Range("D33").Value = '3.71'
$ws.Range("E33").Value = '  +1.09%  '
$ws.Range("D34").Value = '11.00'
$ws.Range("E34").Value = '  -1.71%  '
$ws.Range("D35").Value = '58.32'
$ws.Range("E35").Value = '  +0.59%  '
$ws.Range("E36").Value = '  -2.76%  '
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("D38").Value = '3.590.51'
$ws.Range("E38").Value = '  -2.63%  '
$ws.Range("E39").Value = '  -0.21%  '
$ws.Range("D40").Value = '35.06'
$ws.Range("E40").Value = '  +0.08%  '
$ws.Range("D41").Value = '0.0₃0738'
$ws.Range("E41").Value = '  +5.11%  '
$ws.Range("D42").Value = '3.25'
$ws.Range("E42").Value = '  +0.17%  '
$ws.Range("E43").Value = '  -0.31%  '
$ws.Range("D44").Value = '0.0419'
$ws.Range("E44").Value = '  -0.77%  '
$ws.Range("D45").Value = '0.332'
$ws.Range("E45").Value = '  -1.83%  '
$ws.Range("E46").Value = '  +4.17%  '
$ws.Range("D47").Value = '2.65'
$ws.Range("E47").Value = '  -0.38%  '
$ws.Range("E48").Value = '  -0.47%  '
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  -0.07%  '
$ws.Range("D50").Value = '133.35'
$ws.Range("E50").Value = '  +1.16%  '
$ws.Range("D51").Value = '2.64'
$ws.Range("E51").Value = '  +2.09%  '

# Restore the default (unstyled) appearance for column D so only the
# cell text changes, matching the source workbook which has no explicit
# style index on these data cells.
$ws.Range("D2:D51").Style = "Normal"
